$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1114
$ws1.Range("F3").Value = 0
$ws1.Range("F5").Value = 0
$ws1.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202407/90evzodU1720409956766.jpeg"
$ws1.Range("F7").Value = 8958
$ws1.Range("F8").Value = 231
$ws1.Range("F10").Value = 0
$ws1.Range("F11").Value = 592
$ws1.Range("F12").Value = 0

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0
$ws2.Range("F6").Value = 0

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F5").Value = 14
$ws4.Range("F7").Value = 4900
$ws4.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202407/90evzodU1720409956766.jpeg"
$ws4.Range("F8").Value = 0
$ws4.Range("F10").Value = 8958
$ws4.Range("F11").Value = 231
$ws4.Range("F12").Value = 508
$ws4.Range("F14").Value = 6
$ws4.Range("F15").Value = 2
$ws4.Range("F16").Value = 592
